$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value even when the string looks numeric
# (e.g. "1.00", "0.0625"), mirroring how these values were authored as plain text
# in the source workbook (inline strings), not as numbers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "26.659.90"
$ws.Range("E2").Value = "  +0.88%  "

$ws.Range("D3").Value = "1.642.89"
$ws.Range("E3").Value = "  +1.13%  "

$ws.Range("E4").Value = "  -0.34%  "

Set-TextValue $ws.Range("D5") "214.72"
$ws.Range("E5").Value = "  +0.92%  "

$ws.Range("E6").Value = "  +1.75%  "

Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("E8").Value = "  +0.82%  "

Set-TextValue $ws.Range("D9") "0.0626"
$ws.Range("E9").Value = "  +0.95%  "

Set-TextValue $ws.Range("D10") "19.06"
$ws.Range("E10").Value = "  +0.65%  "

Set-TextValue $ws.Range("D11") "0.0843"
$ws.Range("E11").Value = "  +0.59%  "

$ws.Range("D12").Value = "1.871.15"
$ws.Range("E12").Value = "  +1.04%  "

$ws.Range("D13").Value = "1.645.48"
$ws.Range("E13").Value = "  +1.23%  "

$ws.Range("E14").Value = "  +1.67%  "

$ws.Range("E15").Value = "  +1.59%  "

Set-TextValue $ws.Range("D16") "64.90"
$ws.Range("E16").Value = "  +1.68%  "

$ws.Range("D17").Value = "26.668.25"
$ws.Range("E17").Value = "  +0.89%  "

$ws.Range("E18").Value = "  +0.75%  "

Set-TextValue $ws.Range("D19") "215.21"
$ws.Range("E19").Value = "  -0.01%  "

Set-TextValue $ws.Range("D20") "1.01"
$ws.Range("E20").Value = "  -0.16%  "

Set-TextValue $ws.Range("D21") "4.34"
$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("E23").Value = "  +2.08%  "

$ws.Range("E24").Value = "  +12.94%  "

Set-TextValue $ws.Range("D25") "144.94"
$ws.Range("E25").Value = "  -2.35%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  -0.31%  "

$ws.Range("E28").Value = "  +4.48%  "

$ws.Range("E29").Value = "  +0.88%  "

$ws.Range("E30").Value = "  +1.10%  "

$ws.Range("E31").Value = "  +1.71%  "

Set-TextValue $ws.Range("D32") "3.36"
$ws.Range("E32").Value = "  +1.32%  "

$ws.Range("E33").Value = "  +2.29%  "

$ws.Range("D34").Value = "1.278.56"
$ws.Range("E34").Value = "  +5.13%  "

$ws.Range("E35").Value = "  +2.76%  "

$ws.Range("E36").Value = "  +1.16%  "

$ws.Range("E37").Value = "  +2.96%  "

Set-TextValue $ws.Range("D38") "0.529"
$ws.Range("E38").Value = "  +5.99%  "

Set-TextValue $ws.Range("D39") "0.825"
$ws.Range("E39").Value = "  +4.04%  "

Set-TextValue $ws.Range("D40") "1.01"

Set-TextValue $ws.Range("D41") "0.810"
$ws.Range("E41").Value = "  +2.26%  "

$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("E43").Value = "  +1.26%  "

$ws.Range("D44").Value = "1.781.48"
$ws.Range("E44").Value = "  +1.04%  "

Set-TextValue $ws.Range("D45") "91.65"
$ws.Range("E45").Value = "  -1.22%  "

Set-TextValue $ws.Range("D46") "59.09"
$ws.Range("E46").Value = "  +8.20%  "

$ws.Range("E47").Value = "  +1.44%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0101"
$ws.Range("E48").Value = "  -0.96%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.0516"
$ws.Range("E49").Value = "  +0.95%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "7.73"
$ws.Range("E50").Value = "  +1.90%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.0962"
$ws.Range("E51").Value = "  +1.87%  "
